# Apply updated cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.280.39"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "2.906.43"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "351.89"
$ws.Range("E5").Value = "  -1.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.63"
$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.556"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.84"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  +2.45%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.84"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").Value = "3.370.00"
$ws.Range("E15").Value = "  +3.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.00"
$ws.Range("E16").Value = "  +6.14%  "

$ws.Range("D17").Value = "2.905.71"
$ws.Range("E17").Value = "  +3.30%  "

$ws.Range("D18").Value = "52.336.40"
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.30"
$ws.Range("E20").Value = "  +3.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("E21").Value = "  +3.33%  "

$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.79"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.26"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.76"
$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.69"
$ws.Range("E26").Value = "  +1.94%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.167"
$ws.Range("E28").Value = "  +2.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.58"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").Value = "  +11.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.57"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.55"
$ws.Range("E32").Value = "  +5.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0971"
$ws.Range("E34").Value = "  +10.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.15"
$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0448"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +5.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.68"
$ws.Range("E39").Value = "  -1.11%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("E40").Value = "  +13.24%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.06"
$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.52"
$ws.Range("E42").Value = "  +6.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.63"
$ws.Range("E44").Value = "  +8.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "120.60"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("E47").Value = "  +3.67%  "

$ws.Range("D48").Value = "2.189.33"
$ws.Range("E48").Value = "  +3.69%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.264"
$ws.Range("E49").Value = "  +23.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0337"
$ws.Range("E50").Value = "  +11.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.961"
$ws.Range("E51").Value = "  +1.69%  "

